$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.933.97'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +2.49%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.569.33'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.90%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '600.52'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.73%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '179.25'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.82%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.520'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.65%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.569.55'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.160'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +12.87%  '
$ws.Range('E11').Value = '  -0.12%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.346'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.27%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.03'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.33%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.047.43'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.94%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '26.49'
$ws.Range('D15').Style = 'Normal'
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000182'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +5.43%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '69.884.74'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.85%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.598.61'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.36%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.71'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.59%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.19'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.06%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '365.54'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.47%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.14'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.49%  '
$ws.Range('E23').Value = '  -0.12%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '70.75'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('E26').Value = '  -2.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.29'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.81%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.697.02'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.30%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.67%  '
$ws.Range('B30').Value = 'PEPE'
$ws.Range('C30').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0936'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.21%  '
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '523.19'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.09%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.82'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.15%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.28'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.27%  '
$ws.Range('B34').Value = 'PancakeSwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.81'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.33%  '
$ws.Range('B35').Value = 'FirstDigitalUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.00'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.21%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.121'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.01%  '
$ws.Range('B37').Value = 'Monero'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '162.92'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.84%  '
$ws.Range('B38').Value = 'EthereumClassic'
$ws.Range('C38').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '19.04'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.14%  '
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.92'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.34%  '
$ws.Range('B40').Value = 'ImmutableX'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.36'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.42%  '
$ws.Range('B41').Value = 'Stacks'
$ws.Range('C41').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.77'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.74%  '
$ws.Range('B42').Value = 'USDe'
$ws.Range('C42').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.00'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '4.95'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.56%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.327'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -1.19%  '
$ws.Range('B45').Value = 'dogwifhat'
$ws.Range('C45').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.49'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.14%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.32%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '153.58'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.40%  '
$ws.Range('B48').Value = 'Filecoin'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.64'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +2.18%  '
$ws.Range('B49').Value = 'ARBITRUM'
$ws.Range('C49').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.526'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.65%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₆0260'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.45%  '
$ws.Range('B51').Value = 'Optimism'
$ws.Range('C51').Value = 'https://coinranking.com/coin/n1p-s_gm1+optimism-op'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.63'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +1.67%  '
